# "Add bottom of screen back/next"
# Adds a new `showFooter` setting row to the `settings` sheet, bumps the
# device id in B3, and makes `settings` the active sheet/tab instead of
# `table_specific_translations`.

$wb = $excel.ActiveWorkbook

$settings = $wb.Worksheets.Item("settings")

# Update the survey/form device id stamp.
$settings.Range("B3").Value = 20210221001

# Add a new showFooter=1 setting row right after the existing ones (row 10).
$settings.Range("A10").Value = "showFooter"
$settings.Range("B10").Value = 1

# Make "settings" the active sheet (was "table_specific_translations") and
# move the selection/active cell there.
$settings.Activate()
$settings.Range("B11").Select()
